# cierre 18 Jun 22
# Advance the payroll week label by one week, and update the two
# manually-entered pay figures for the new week. Formula-driven cells
# (the mirrored "SEMANA..." labels, the TODAY()-based dates, and the
# SUM() totals) recalculate on their own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Week banner: week 23 (06-12 Jun 2022) -> week 24 (13-19 Jun 2022).
# B9 is the master cell; H9, B27, H27 and B43 all reference it (or each
# other) via formulas, so they update automatically on recalculation.
$ws.Range("B9").Value = "SEMANA   24  DEL    13      Al   19   DE   JUNIO          2022"

# Manually entered pay amounts for the new week.
$ws.Range("K21").Value = 980
$ws.Range("E40").Value = 0

# Restore the view state (scroll position / active selection) as left by
# the user after editing.
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 25
$ws.Range("I48:I49").Select()
$excel.ActiveCell = $ws.Range("I49")

$wb.Application.Calculate()
